$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to rounded (2 decimal place) figures
$ws.Range("B5").Value = 6.86
$ws.Range("C5").Value = 5.25
$ws.Range("D5").Value = 0.4
$ws.Range("F5").Value = 12.35
$ws.Range("G5").Value = 5.15
$ws.Range("H5").Value = 21.01
$ws.Range("J5").Value = 3.64
$ws.Range("K5").Value = 5.63
$ws.Range("L5").Value = 5.95
$ws.Range("M5").Value = 6.14
$ws.Range("N5").Value = 1.66
$ws.Range("O5").Value = 5.15
$ws.Range("P5").Value = 7.49
$ws.Range("Q5").Value = 4.33
$ws.Range("R5").Value = 0.36
$ws.Range("S5").Value = 0.07000000000000001
$ws.Range("T5").Value = 72.61
$ws.Range("U5").Value = 14.52
$ws.Range("V5").Value = 4.7
$ws.Range("W5").Value = 9.789999999999999
$ws.Range("X5").Value = 5.42
$ws.Range("Y5").Value = 0.6899999999999999
$ws.Range("AA5").Value = 3.97
$ws.Range("AC5").Value = 4.25
$ws.Range("AF5").Value = 18.8
$ws.Range("AH5").Value = 5.91

# Remove the now-obsolete last data row (row 6); dimension shrinks to A1:AH5 automatically
$ws.Rows.Item(6).Delete()

# Narrow column T (20th column) from width 8 to width 7
$ws.Columns.Item(20).ColumnWidth = 6.17
